# Strike through the "Mettre une sortie..." to-do item: it (and related
# items mentioned in the commit message - timer, survivor counter, safe
# zone, etc.) has been completed, so mark it as done with strikethrough
# formatting on both the paragraph mark and its run.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("Mettre une sortie, qui", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $para = $find.Parent.Paragraphs(1)
    $para.Range.Font.StrikeThrough = 1
}
